$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to "All Published Values" ----------------------
$ws = $wb.Worksheets.Item("All Published Values")
$newRow = 16

# Force the row to Text format before writing so Excel doesn't silently
# reinterpret the date-/number-looking strings as a real date serial or a
# double (the source data is stored as literal text, like the rest of the
# sheet). The style is put back to "Normal" afterwards so the new cells end
# up on the same (default) style as every other data row.
$ws.Range("A16:J16").NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value  = "2026-01-02"
$ws.Cells.Item($newRow, 2).Value  = "2026-01-02 22:16:58"
$ws.Cells.Item($newRow, 3).Value  = "697.85"
$ws.Cells.Item($newRow, 4).Value  = "697.85"
$ws.Cells.Item($newRow, 5).Value  = "700.79"
$ws.Cells.Item($newRow, 6).Value  = "700.79"
$ws.Cells.Item($newRow, 7).Value  = "702.88"
$ws.Cells.Item($newRow, 8).Value  = "2026/01/02 22:16:58"
$ws.Cells.Item($newRow, 9).Value  = "2026-01-02 14:19:32"
$ws.Cells.Item($newRow, 10).Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

$ws.Range("A16:J16").Style = "Normal"

# --- 2. Re-stretch the AutoFilter range to cover the new row --------------
# Toggle off first: calling Range.AutoFilter() a second time on a sheet that
# already has an active AutoFilter flips it OFF (VBA toggle semantics), so
# clear it, then re-apply across the full, now-16-row range.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:J16").AutoFilter()

# --- 3. Keep the hidden _FilterDatabase defined name in sync --------------
$nm = $wb.Names.Item("All Published Values!_FilterDatabase")
$nm.RefersTo = "='All Published Values'!`$A`$1:`$J`$16"

# --- 4. Update the "Daily Summary" publishes count -------------------------
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Cells.Item(4, 2).Value = 15
